$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to keep text formatting so numeric-looking price strings
# (e.g. "1.004") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '28.037.06'
$ws.Range("E2").Value = '  -0.17%  '

# Row 3
$ws.Range("D3").Value = '1.872.49'
$ws.Range("E3").Value = '  -0.10%  '

# Row 4
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("E5").Value = '  -0.42%  '

# Row 6
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Range("D7").Value = '0.5152'
$ws.Range("E7").Value = '  +2.19%  '

# Row 8
$ws.Range("D8").Value = '0.3850'
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").Value = '0.08297'
$ws.Range("E9").Value = '  -3.67%  '

# Row 10
$ws.Range("D10").Value = '1.111'
$ws.Range("E10").Value = '  -0.31%  '

# Row 11
$ws.Range("D11").Value = '41.53'
$ws.Range("E11").Value = '  +0.14%  '

# Row 12
$ws.Range("E12").Value = '  -1.60%  '

# Row 13
$ws.Range("E13").Value = '  -0.47%  '

# Row 14
$ws.Range("D14").Value = '1.873.81'
$ws.Range("E14").Value = '  -0.17%  '

# Row 15
$ws.Range("D15").Value = '7.304'
$ws.Range("E15").Value = '  +1.47%  '

# Row 16
$ws.Range("E16").Value = '  -0.01%  '

# Row 17
$ws.Range("E17").Value = '  -0.11%  '

# Row 18
$ws.Range("D18").Value = '90.85'
$ws.Range("E18").Value = '  -0.14%  '

# Row 19
$ws.Range("D19").Value = '0.06652'
$ws.Range("E19").Value = '  +0.34%  '

# Row 20
$ws.Range("D20").Value = '17.73'

# Row 21
$ws.Range("E21").Value = '  -0.01%  '

# Row 22
$ws.Range("D22").Value = '6.035'
$ws.Range("E22").Value = '  -0.90%  '

# Row 23
$ws.Range("D23").Value = '28.080.59'
$ws.Range("E23").Value = '  -0.14%  '

# Row 24
$ws.Range("D24").Value = '11.09'
$ws.Range("E24").Value = '  -2.73%  '

# Row 25
$ws.Range("D25").Value = '2.253'
$ws.Range("E25").Value = '  -0.55%  '

# Row 26
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.075.80'
$ws.Range("E26").Value = '  -0.92%  '

# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.514'
$ws.Range("E27").Value = '  -2.88%  '

# Row 28
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '157.38'
$ws.Range("E28").Value = '  +0.14%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '20.55'
$ws.Range("E29").Value = '  -0.89%  '

# Row 30
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '125.01'
$ws.Range("E30").Value = '  -0.91%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.1065'
$ws.Range("E31").Value = '  +1.04%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '1.034'
$ws.Range("E32").Value = '  -2.50%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.752'
$ws.Range("E33").Value = '  +2.64%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.591'
$ws.Range("E34").Value = '  -0.13%  '

# Row 35
$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").Value = '9.478'
$ws.Range("E35").Value = '  -1.52%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02424'
$ws.Range("E36").Value = '  -0.74%  '

# Row 37
$ws.Range("D37").Value = '0.06525'
$ws.Range("E37").Value = '  -0.86%  '

# Row 38
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = '0.2202'
$ws.Range("E38").Value = '  +1.07%  '

# Row 39
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.6571'
$ws.Range("E39").Value = '  +2.94%  '

# Row 40
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '1.203'
$ws.Range("E40").Value = '  -0.58%  '

# Row 41
$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").Value = '5.017'
$ws.Range("E41").Value = '  +2.51%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.212'
$ws.Range("E42").Value = '  -2.46%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '11.23'
$ws.Range("E43").Value = '  -2.36%  '

# Row 44
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = '0.6135'
$ws.Range("E44").Value = '  +2.21%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.10'
$ws.Range("E45").Value = '  -0.45%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.283'
$ws.Range("E46").Value = '  +0.11%  '

# Row 47
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.676'
$ws.Range("E47").Value = '  +0.07%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.026'
$ws.Range("E48").Value = '  +1.87%  '

# Row 49
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.218'
$ws.Range("E49").Value = '  -0.56%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '121.12'
$ws.Range("E50").Value = '  -0.49%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '78.48'
$ws.Range("E51").Value = '  -2.32%  '

